# Apply updated currentAveragePrice / LeveProfit figures across several
# crafting-sheet tables (scheduled-runner refresh of market-board data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (item id 38956)
$ws.Range("H17").Value = 2150
$ws.Range("J17").Value = 2150
$ws.Range("L17").Value = 6450
$ws.Range("N17").Value = -6786

# Row 33 (item id 5512)
$ws.Range("H33").Value = 121.4
$ws.Range("I33").Value = 122.666664
$ws.Range("K33").Value = 122.666664
$ws.Range("M33").Value = 106.333336

# Row 64 (item id 5506)
$ws.Range("H64").Value = 2993.3333
$ws.Range("I64").Value = 2993.3333
$ws.Range("K64").Value = 2993.3333
$ws.Range("M64").Value = -2745.3333

# Row 67 (item id 5506)
$ws.Range("H67").Value = 2993.3333
$ws.Range("I67").Value = 2993.3333
$ws.Range("K67").Value = 2993.3333
$ws.Range("M67").Value = -2135.3333

# Row 74 (item id 5507)
$ws.Range("H74").Value = 4699.25
$ws.Range("I74").Value = 4599.3335
$ws.Range("J74").Value = 4999
$ws.Range("K74").Value = 4599.3335
$ws.Range("L74").Value = 4999
$ws.Range("M74").Value = -3663.3335
$ws.Range("N74").Value = -6871

# Row 76 (item id 12602)
$ws.Range("H76").Value = 7285.5713
$ws.Range("I76").Value = 8000
$ws.Range("K76").Value = 8000
$ws.Range("M76").Value = -7685

# Row 77 (item id 5507)
$ws.Range("H77").Value = 4699.25
$ws.Range("I77").Value = 4599.3335
$ws.Range("J77").Value = 4999
$ws.Range("K77").Value = 22996.6675
$ws.Range("L77").Value = 24995
$ws.Range("M77").Value = -18316.6675
$ws.Range("N77").Value = -34355

# Row 79 (item id 12602)
$ws.Range("H79").Value = 7285.5713
$ws.Range("I79").Value = 8000
$ws.Range("K79").Value = 8000
$ws.Range("M79").Value = -6908

# Row 137 (item id 44013)
$ws.Range("H137").Value = 3167.6667
$ws.Range("I137").Value = 5333.3335
$ws.Range("J137").Value = 1002
$ws.Range("K137").Value = 16000.0005
$ws.Range("L137").Value = 3006
$ws.Range("M137").Value = -13450.0005
$ws.Range("N137").Value = -8106

# Row 138 (item id 44169)
$ws.Range("H138").Value = 3849541.2
$ws.Range("J138").Value = 3566.45
$ws.Range("L138").Value = 10699.35
$ws.Range("N138").Value = -20979.35

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id 44147)
$ws.Range("H32").Value = 3556.8
$ws.Range("I32").Value = 3556.8
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3556.8
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3269.8
$ws.Range("N32").Value = $null

# Row 46 (item id 3498)
$ws.Range("H46").Value = 4569
$ws.Range("I46").Value = 4569
$ws.Range("K46").Value = 4569
$ws.Range("M46").Value = -4250

# Row 63 (item id 12528)
$ws.Range("H63").Value = 5199.857
$ws.Range("I63").Value = 5233.1665
$ws.Range("K63").Value = 5233.1665
$ws.Range("M63").Value = -4547.1665

# Row 66 (item id 12528)
$ws.Range("H66").Value = 5199.857
$ws.Range("I66").Value = 5233.1665
$ws.Range("K66").Value = 26165.8325
$ws.Range("M66").Value = -22733.8325

# Row 122 (item id 36168)
$ws.Range("H122").Value = 2616.5
$ws.Range("I122").Value = 2403.8333
$ws.Range("J122").Value = 3254.5
$ws.Range("K122").Value = 7211.499899999999
$ws.Range("L122").Value = 9763.5
$ws.Range("M122").Value = -4761.499899999999
$ws.Range("N122").Value = -14663.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (item id 19947)
$ws.Range("H105").Value = 2031.5714
$ws.Range("I105").Value = 1953.5
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1953.5
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -206.5
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (item id 44021)
$ws.Range("H58").Value = 7210.636
$ws.Range("I58").Value = 4760.4287
$ws.Range("J58").Value = 11498.5
$ws.Range("K58").Value = 4760.4287
$ws.Range("L58").Value = 11498.5
$ws.Range("M58").Value = -4557.4287
$ws.Range("N58").Value = -11904.5

# Row 109 (item id 27203)
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

# Row 134 (item id 44020)
$ws.Range("H134").Value = 3802.7856
$ws.Range("I134").Value = 1425.1111
$ws.Range("J134").Value = 8082.6
$ws.Range("K134").Value = 4275.3333
$ws.Range("L134").Value = 24247.8
$ws.Range("M134").Value = -1740.3333
$ws.Range("N134").Value = -29317.8

# Row 136 (item id 44021)
$ws.Range("H136").Value = 7210.636
$ws.Range("I136").Value = 4760.4287
$ws.Range("J136").Value = 11498.5
$ws.Range("K136").Value = 14281.2861
$ws.Range("L136").Value = 34495.5
$ws.Range("M136").Value = -11731.2861
$ws.Range("N136").Value = -39595.5

$ws = $wb.Worksheets.Item("GSM")
# Row 19 (item id 2668)
$ws.Range("H19").Value = 3425.6155
$ws.Range("I19").Value = 3653.4
$ws.Range("K19").Value = 3653.4
$ws.Range("M19").Value = -3365.4

# Row 102 (item id 36169)
$ws.Range("H102").Value = 3289.3333
$ws.Range("I102").Value = 1903.125
$ws.Range("K102").Value = 1903.125
$ws.Range("M102").Value = -281.125

$ws = $wb.Worksheets.Item("LTW")
# Row 11 (item id 3542)
$ws.Range("H11").Value = 524.5
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 999
$ws.Range("K11").Value = 50
$ws.Range("L11").Value = 999
$ws.Range("M11").Value = 90
$ws.Range("N11").Value = -1279

# Row 68 (item id 12563)
$ws.Range("H68").Value = 2821.2856
$ws.Range("I68").Value = 2958.1667
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2958.1667
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -2209.1667
$ws.Range("N68").Value = -3498

# Row 71 (item id 12563)
$ws.Range("H71").Value = 2821.2856
$ws.Range("I71").Value = 2958.1667
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 14790.8335
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -11046.8335
$ws.Range("N71").Value = -17488

# Row 80 (item id 12027)
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null

# Row 82 (item id 12565)
$ws.Range("H82").Value = 1350
$ws.Range("I82").Value = 1200
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 1200
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -839
$ws.Range("N82").Value = -2222

# Row 83 (item id 12027)
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null

# Row 85 (item id 12565)
$ws.Range("H85").Value = 1350
$ws.Range("I85").Value = 1200
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 1200
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = 48
$ws.Range("N85").Value = -3996

# Row 99 (item id 19636)
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = $null

# Row 122 (item id 36247)
$ws.Range("H122").Value = 3549.875
$ws.Range("I122").Value = 3399.8333
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10199.4999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7749.499899999999
$ws.Range("N122").Value = -16900

# Row 132 (item id 44058)
$ws.Range("H132").Value = 7013.467
$ws.Range("I132").Value = 4933.0835
$ws.Range("K132").Value = 14799.2505
$ws.Range("M132").Value = -12269.2505

$ws = $wb.Worksheets.Item("WVR")
# Row 12 (item id 3316)
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = $null

# Row 52 (item id 2816)
$ws.Range("H52").Value = 16759.2
$ws.Range("I52").Value = 16759.2
$ws.Range("K52").Value = 16759.2
$ws.Range("M52").Value = -16533.2
